# Uniswap workbook update (10.11.24 -> 26.12.24 gUSD data)
# 1) Rename the "gUSD" sheet to "gUSD 26.12.24"
# 2) Append 8 new days of data (rows 33-40) to the gUSD sheet, following
#    the exact same column layout / formula pattern already used in the
#    sheet (A: date = prior date + 1; B-F: raw inputs; H: B-Bprev;
#    I: C-Cprev; K: B+C). Row 40 only carries the date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD")
$ws.Name = "gUSD 26.12.24"

# New data rows 33-39 (date is derived with =A(n-1)+1 just like the rest
# of the sheet; H/I are only computed from row 33 onward, matching the
# existing pattern where H/I start one row after B/C/D/E/F begin).
$newRows = @(
    @{ Row = 33; B = 6.04;              C = 40.36; D = 14.64; E = 12.71; F = 8.61 }
    @{ Row = 34; B = 5.22;              C = 41.64; D = 14.84; E = 22.58; F = 10.1 }
    @{ Row = 35; B = 4.4000000000000004; C = 43.08; D = 14.87; E = 26.01; F = 12.98 }
    @{ Row = 36; B = 3.6;               C = 43.35; D = 15.09; E = 4.7;   F = 12.72 }
    @{ Row = 37; B = 2.66;              C = 43.65; D = 15.09; E = 4.8;   F = 12.72 }
    @{ Row = 38; B = 1.76;              C = 44.05; D = 15.45; E = 6.59;  F = 12.23 }
    @{ Row = 39; B = 0.83;              C = 44.97; D = 15.54; E = 15.49; F = 12.89 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prev = $row - 1

    $ws.Cells.Item($row, 1).Formula = "=A$prev+1"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 8).Formula = "=B$row-B$prev"
    $ws.Cells.Item($row, 9).Formula = "=C$row-C$prev"
    $ws.Cells.Item($row, 11).Formula = "=B$row+C$row"
}

# Row 40 only gets the date column filled in (mirrors the source diff).
$ws.Cells.Item(40, 1).Formula = "=A39+1"

# Restore the view: scrolled down with A40 selected (best effort - the
# scroll position itself is cosmetic UI state).
$ws.Activate()
$ws.Range("A40").Select()
$excel.ActiveWindow.ScrollRow = 19
